$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell A5 used to hold " lilis88@mail.ru " (stray leading space, no link).
# Clean up the text to "lilis88@mail.ru " (trailing space only) and turn
# it into a live mailto: hyperlink, same as the other email cells above it.
$ws.Range("A5").Value = "lilis88@mail.ru "
$null = $ws.Hyperlinks.Add($ws.Range("A5"), "mailto:lilis88@mail.ru")

# Leave the active cell where the edit session ended up.
$null = $ws.Range("A7").Select()
